$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179.75
$ws.Range("I2").Value = 120
$ws.Range("J2").Value = 199.66667
$ws.Range("K2").Value = 120
$ws.Range("L2").Value = 199.66667
$ws.Range("M2").Value = -7
$ws.Range("N2").Value = -425.66667

$ws.Range("H107").Value = 509.25
$ws.Range("I107").Value = 532.16
$ws.Range("K107").Value = 532.16
$ws.Range("M107").Value = 1387.84

$ws.Range("H125").Value = 592.3125
$ws.Range("I125").Value = 392.55554
$ws.Range("J125").Value = 849.1429000000001
$ws.Range("K125").Value = 3532.99986
$ws.Range("L125").Value = 7642.2861
$ws.Range("M125").Value = -1072.99986
$ws.Range("N125").Value = -12562.2861

$ws.Range("H135").Value = 560.5
$ws.Range("I135").Value = 560.5
$ws.Range("K135").Value = 5044.5
$ws.Range("M135").Value = -2509.5

$ws.Range("H137").Value = 40001540
$ws.Range("I137").Value = 52632504
$ws.Range("J137").Value = 3490.6667
$ws.Range("K137").Value = 157897512
$ws.Range("L137").Value = 10472.0001
$ws.Range("M137").Value = -157894962
$ws.Range("N137").Value = -15572.0001

$ws.Range("H138").Value = 3135.6792
$ws.Range("I138").Value = 1304.909
$ws.Range("J138").Value = 4434.9355
$ws.Range("K138").Value = 3914.727
$ws.Range("L138").Value = 13304.8065
$ws.Range("M138").Value = 1225.273
$ws.Range("N138").Value = -23584.8065

$ws.Range("H141").Value = 731.6667
$ws.Range("I141").Value = 95
$ws.Range("K141").Value = 285
$ws.Range("M141").Value = 4895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8623123
$ws.Range("I32").Value = 10418056
$ws.Range("K32").Value = 10418056
$ws.Range("M32").Value = -10417769

$ws.Range("H61").Value = 3079.946
$ws.Range("I61").Value = 2802.56
$ws.Range("K61").Value = 2802.56
$ws.Range("M61").Value = -2590.56

$ws.Range("H102").Value = 21115.047
$ws.Range("I102").Value = 23601.445
$ws.Range("K102").Value = 23601.445
$ws.Range("M102").Value = -21979.445

$ws.Range("H136").Value = 3079.946
$ws.Range("I136").Value = 2802.56
$ws.Range("K136").Value = 8407.68
$ws.Range("M136").Value = -5857.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 42364.273
$ws.Range("I82").Value = 3158.1667
$ws.Range("J82").Value = 89411.60000000001
$ws.Range("K82").Value = 3158.1667
$ws.Range("L82").Value = 89411.60000000001
$ws.Range("M82").Value = -2775.1667
$ws.Range("N82").Value = -90177.60000000001

$ws.Range("H85").Value = 42364.273
$ws.Range("I85").Value = 3158.1667
$ws.Range("J85").Value = 89411.60000000001
$ws.Range("K85").Value = 3158.1667
$ws.Range("L85").Value = 89411.60000000001
$ws.Range("M85").Value = -1832.1667
$ws.Range("N85").Value = -92063.60000000001

$ws.Range("H105").Value = 1697.2727
$ws.Range("I105").Value = 1636.7222
$ws.Range("K105").Value = 1636.7222
$ws.Range("M105").Value = 110.2778000000001

$ws.Range("H132").Value = 111387.5
$ws.Range("J132").Value = 111387.5
$ws.Range("L132").Value = 111387.5
$ws.Range("N132").Value = -121507.5

$ws.Range("H134").Value = 25996.56
$ws.Range("I134").Value = 32957.195
$ws.Range("J134").Value = 4418.6
$ws.Range("K134").Value = 98871.58499999999
$ws.Range("L134").Value = 13255.8
$ws.Range("M134").Value = -96336.58499999999
$ws.Range("N134").Value = -18325.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 65793.92
$ws.Range("I31").Value = 83863.30499999999
$ws.Range("K31").Value = 83863.30499999999
$ws.Range("M31").Value = -83568.30499999999

$ws.Range("H34").Value = 65793.92
$ws.Range("I34").Value = 83863.30499999999
$ws.Range("K34").Value = 83863.30499999999
$ws.Range("M34").Value = -83661.30499999999

$ws.Range("H58").Value = 2258.0286
$ws.Range("I58").Value = 1352.7142
$ws.Range("J58").Value = 3616
$ws.Range("K58").Value = 1352.7142
$ws.Range("L58").Value = 3616
$ws.Range("M58").Value = -1149.7142
$ws.Range("N58").Value = -4022

$ws.Range("H88").Value = 19460.75
$ws.Range("J88").Value = 19460.75
$ws.Range("L88").Value = 19460.75
$ws.Range("N88").Value = -20272.75

$ws.Range("H91").Value = 19460.75
$ws.Range("J91").Value = 19460.75
$ws.Range("L91").Value = 19460.75
$ws.Range("N91").Value = -22268.75

$ws.Range("H122").Value = 1916.3846
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H132").Value = 44644024
$ws.Range("I132").Value = 38462676
$ws.Range("K132").Value = 115388028
$ws.Range("M132").Value = -115385498

$ws.Range("H136").Value = 2258.0286
$ws.Range("I136").Value = 1352.7142
$ws.Range("J136").Value = 3616
$ws.Range("K136").Value = 4058.1426
$ws.Range("L136").Value = 10848
$ws.Range("M136").Value = -1508.1426
$ws.Range("N136").Value = -15948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1730.1111
$ws.Range("I70").Value = 529.3333
$ws.Range("K70").Value = 1587.9999
$ws.Range("M70").Value = -1272.9999

$ws.Range("H73").Value = 1730.1111
$ws.Range("I73").Value = 529.3333
$ws.Range("K73").Value = 1587.9999
$ws.Range("M73").Value = -495.9999

$ws.Range("H75").Value = 8382.571
$ws.Range("I75").Value = 6725
$ws.Range("K75").Value = 20175
$ws.Range("M75").Value = -19177

$ws.Range("H78").Value = 8382.571
$ws.Range("I78").Value = 6725
$ws.Range("K78").Value = 60525
$ws.Range("M78").Value = -55533

$ws.Range("H107").Value = 301.53845
$ws.Range("J107").Value = 368.3
$ws.Range("L107").Value = 1104.9
$ws.Range("N107").Value = -4944.9

$ws.Range("H118").Value = 1195.6666
$ws.Range("I118").Value = 1194.8
$ws.Range("K118").Value = 3584.4
$ws.Range("M118").Value = -2341.4

$ws.Range("H129").Value = 1047.3158
$ws.Range("I129").Value = 757.86664
$ws.Range("J129").Value = 2132.75
$ws.Range("K129").Value = 2273.59992
$ws.Range("L129").Value = 6398.25
$ws.Range("M129").Value = 2726.40008
$ws.Range("N129").Value = -16398.25

$ws.Range("H132").Value = 906
$ws.Range("I132").Value = 907.3333
$ws.Range("K132").Value = 8165.9997
$ws.Range("M132").Value = -5635.9997

$ws.Range("H139").Value = 2190
$ws.Range("I139").Value = 1744
$ws.Range("K139").Value = 5232
$ws.Range("M139").Value = -92

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 61187.35
$ws.Range("I122").Value = 62279.41
$ws.Range("K122").Value = 186838.23
$ws.Range("M122").Value = -184388.23

$ws.Range("H132").Value = 29428590
$ws.Range("I132").Value = 45469540
$ws.Range("K132").Value = 136408620
$ws.Range("M132").Value = -136406090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 69552.664
$ws.Range("I122").Value = 2238.6
$ws.Range("J122").Value = 204180.8
$ws.Range("K122").Value = 6715.799999999999
$ws.Range("L122").Value = 612542.3999999999
$ws.Range("M122").Value = -4265.799999999999
$ws.Range("N122").Value = -617442.3999999999

$ws.Range("H132").Value = 2696.1282
$ws.Range("I132").Value = 2613.875
$ws.Range("K132").Value = 7841.625
$ws.Range("M132").Value = -5311.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1481.6666
$ws.Range("I122").Value = 1183.8667
$ws.Range("K122").Value = 3551.6001
$ws.Range("M122").Value = -1101.6001

$ws.Range("H123").Value = 34500
$ws.Range("J123").Value = 34500
$ws.Range("L123").Value = 34500
$ws.Range("N123").Value = -44300

$ws.Range("H132").Value = 5130985.5
$ws.Range("I132").Value = 7409777
$ws.Range("J132").Value = 3704.6667
$ws.Range("K132").Value = 22229331
$ws.Range("L132").Value = 11114.0001
$ws.Range("M132").Value = -22226801
$ws.Range("N132").Value = -16174.0001
